$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of rows 2 & 3 with rows 4 & 5 for columns
# D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg).
# All other columns (A, B, C, E, F, G, H, I, J, K, Q, T) stay the same.

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $rangeTop = $ws.Range($col + "2")
    $rangeBottom = $ws.Range($col + "4")
    $valTop = $rangeTop.Value()
    $valBottom = $rangeBottom.Value()
    $rangeTop.Value = $valBottom
    $rangeBottom.Value = $valTop

    $rangeTop2 = $ws.Range($col + "3")
    $rangeBottom2 = $ws.Range($col + "5")
    $valTop2 = $rangeTop2.Value()
    $valBottom2 = $rangeBottom2.Value()
    $rangeTop2.Value = $valBottom2
    $rangeBottom2.Value = $valTop2
}

$wb.Save()
